$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New column BB: header date (BB1) ---
# Set the value first, then copy only the number-format/style from BA1 so
# Excel doesn't invent a brand-new style for the literal date serial.
$ws.Range("BB1").Value = 45986
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)  # xlPasteFormats

# --- New column BB rows 2-81: same values as column BA (unstyled numeric cells) ---
$ws.Range("BA2:BA81").Copy()
$ws.Range("BB2:BB81").PasteSpecial(-4163)  # xlPasteValues (no style carried over)

# --- Rows 82 & 83: BB differs from BA ---
$ws.Range("BB82").Value = 1.538981993999982
$ws.Range("BB83").Value = 1.68501852020853

# --- New row 84: only A84 (date) and BB84 (value) are populated ---
$ws.Range("A84").Value = 45884
$ws.Range("A83").Copy()
$ws.Range("A84").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("BB84").Value = 0.03331000006224372

$excel.CutCopyMode = 0
